# Apply text replacements to ActasHechos1.docx per commit diff.
# "cambios en identificacion-en cola/urgente/showbyfolio"

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

# Time of the act
Replace-Text "12:28:" "11:54:"

# Date of the act
Replace-Text "lunes 23 de abril del año 2018" "martes 24 de abril del año 2018"

# Name of the declarant (appears multiple times -> ReplaceAll handles all)
Replace-Text "CESAR SERRANO CARRION" "PAOLA SUAREZ BUENO"

# Identification type
Replace-Text "CREDENCIAL DE ELECTOR " "INE "

# Folio number
Replace-Text "AISDJGMSLFAOISGF84651" "468473216546"

# Issuer of the ID
Replace-Text "EL INSTITUTO FEDERAL ELECTORAL" "MI"

# Age (standalone run "28" years old)
Replace-Text "28 años de edad" "24 años de edad"

# Date of birth
Replace-Text "24 de agosto del año 1989" "23 de febrero del año 1994"

# Street address
Replace-Text "JOSE MANCISIDOR 1" "SAKDASJDASDJASKJD 12 interior 12"

# Colonia
Replace-Text "ISLETA, C.P." "NIÑOS HEROES, C.P."

# Postal code
Replace-Text "91090" "91015"

# Occupation
Replace-Text "ANALISTA DE SISTEMAS INFORMÁTICOS" "ABOGADO CIVILISTA"

# Marital status
Replace-Text "CASADO" "CONCUBINATO"

# Education level
Replace-Text "SECUNDARIA INCOMPLETA" "POSGRADO COMPLETO"

# Phone number
Replace-Text "2281909090" "1651432165848"

# Declaration text
Replace-Text "ME ROBARON MI CARTERA JUNTO CON EL INE" "XKCJA,SLCKLA ALSJCAKLNSCÑLASMDNALSCAHSCLAKSDLOASJDKJASLKD"
